$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.215.86'
$ws.Range("E2").Value = '  -6.09%  '
$ws.Range("D3").Value = '2.214.18'
$ws.Range("E3").Value = '  -6.35%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = "'242.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("D6").Value = "'0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.62%  '
$ws.Range("D7").Value = "'69.91"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.97%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -8.05%  '
$ws.Range("D10").Value = "'38.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.27%  '
$ws.Range("D11").Value = "'0.0952"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.97%  '
$ws.Range("D12").Value = "'57.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.85%  '
$ws.Range("E13").Value = '  -3.95%  '
$ws.Range("D14").Value = "'6.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.73%  '
$ws.Range("D15").Value = '2.542.65'
$ws.Range("E15").Value = '  -6.43%  '
$ws.Range("D16").Value = "'14.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -9.55%  '
$ws.Range("D17").Value = "'0.840"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -9.52%  '
$ws.Range("D18").Value = '2.215.97'
$ws.Range("E18").Value = '  -8.29%  '
$ws.Range("D19").Value = '41.165.01'
$ws.Range("E19").Value = '  -6.10%  '
$ws.Range("D20").Value = '0.0₃0950'
$ws.Range("E20").Value = '  -8.22%  '
$ws.Range("D21").Value = "'72.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.95%  '
$ws.Range("E22").Value = '  -7.93%  '
$ws.Range("D23").Value = "'231.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.77%  '
$ws.Range("D24").Value = "'2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.44%  '
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").Value = "'3.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.24%  '
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D28").Value = "'9.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.54%  '
$ws.Range("E29").Value = '  -4.95%  '
$ws.Range("D30").Value = "'172.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.66%  '
$ws.Range("D31").Value = "'20.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.80%  '
$ws.Range("E32").Value = '  -7.80%  '
$ws.Range("D33").Value = "'0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.77%  '
$ws.Range("D34").Value = "'0.0710"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.55%  '
$ws.Range("D35").Value = "'5.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.93%  '
$ws.Range("D36").Value = "'4.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.77%  '
$ws.Range("E37").Value = '  +2.84%  '
$ws.Range("D38").Value = "'23.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +16.09%  '
$ws.Range("D39").Value = "'0.0279"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.57%  '
$ws.Range("E40").Value = '  -5.52%  '
$ws.Range("E41").Value = '  -11.86%  '
$ws.Range("D42").Value = "'64.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.99%  '
$ws.Range("D43").Value = "'4.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -11.78%  '
$ws.Range("D44").Value = "'0.197"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.60%  '
$ws.Range("D45").Value = "'8.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.84%  '
$ws.Range("E46").Value = '  -7.02%  '
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("E48").Value = '  +10.60%  '
$ws.Range("D49").Value = "'4.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.99%  '
$ws.Range("E50").Value = '  -5.94%  '
$ws.Range("E51").Value = '  -5.64%  '
